# Auto-generated script applying the Golem_Profits market-data refresh
# (scheduled runner updates currentAveragePrice / Leve price / profit columns)
$wb = $excel.ActiveWorkbook

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 27.8
$ws.Range("I8").Value = 9.75
$ws.Range("K8").Value = 29.25
$ws.Range("M8").Value = 109.75

# ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 338.18182
$ws.Range("I42").Value = 218.125
$ws.Range("J42").Value = 658.3333
$ws.Range("K42").Value = 654.375
$ws.Range("L42").Value = 1974.9999
$ws.Range("M42").Value = -424.375
$ws.Range("N42").Value = -2434.9999

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 57784.668
$ws.Range("I87").Value = 35000
$ws.Range("K87").Value = 35000
$ws.Range("M87").Value = -33752

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 57784.668
$ws.Range("I90").Value = 35000
$ws.Range("K90").Value = 105000
$ws.Range("M90").Value = -98760

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 43133.285
$ws.Range("J107").Value = 419.8
$ws.Range("L107").Value = 419.8
$ws.Range("N107").Value = -4259.8

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5361.8
$ws.Range("I116").Value = 5361.8
$ws.Range("K116").Value = 5361.8
$ws.Range("M116").Value = -1919.8

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5011.294
$ws.Range("I32").Value = 5011.294
$ws.Range("K32").Value = 5011.294
$ws.Range("M32").Value = -4724.294

# ARM row 48
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H48").Value = 199999
$ws.Range("J48").Value = 199999
$ws.Range("L48").Value = 199999
$ws.Range("N48").Value = -200767

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1482.4286
$ws.Range("I63").Value = 1680
$ws.Range("J63").Value = 988.5
$ws.Range("K63").Value = 1680
$ws.Range("L63").Value = 988.5
$ws.Range("M63").Value = -994
$ws.Range("N63").Value = -2360.5

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1482.4286
$ws.Range("I66").Value = 1680
$ws.Range("J66").Value = 988.5
$ws.Range("K66").Value = 8400
$ws.Range("L66").Value = 4942.5
$ws.Range("M66").Value = -4968
$ws.Range("N66").Value = -11806.5

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 125000500
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1070.3334
$ws.Range("I122").Value = 1070.3334
$ws.Range("K122").Value = 3211.0002
$ws.Range("M122").Value = -761.0001999999999

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -826

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26108.066
$ws.Range("I82").Value = 14909.385
$ws.Range("K82").Value = 14909.385
$ws.Range("M82").Value = -14526.385

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 26108.066
$ws.Range("I85").Value = 14909.385
$ws.Range("K85").Value = 14909.385
$ws.Range("M85").Value = -13583.385

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 66651.766
$ws.Range("J94").Value = 3083.3333
$ws.Range("L94").Value = 3083.3333
$ws.Range("N94").Value = -3985.3333

# BSM row 97
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 10669.2
$ws.Range("I97").Value = 10669.2
$ws.Range("K97").Value = 10669.2
$ws.Range("M97").Value = -9678.200000000001

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2332.5
$ws.Range("I134").Value = 1999
$ws.Range("K134").Value = 5997
$ws.Range("M134").Value = -3462

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 17.714285
$ws.Range("I7").Value = 17
$ws.Range("J7").Value = 19.5
$ws.Range("K7").Value = 51
$ws.Range("L7").Value = 58.5
$ws.Range("M7").Value = 61
$ws.Range("N7").Value = -282.5

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 407.85715
$ws.Range("I86").Value = 250
$ws.Range("J86").Value = 471
$ws.Range("K86").Value = 750
$ws.Range("L86").Value = 1413
$ws.Range("M86").Value = 436
$ws.Range("N86").Value = -3785

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 407.85715
$ws.Range("I89").Value = 250
$ws.Range("J89").Value = 471
$ws.Range("K89").Value = 2250
$ws.Range("L89").Value = 4239
$ws.Range("M89").Value = 3678
$ws.Range("N89").Value = -16095

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 490
$ws.Range("I122").Value = 490
$ws.Range("K122").Value = 4410
$ws.Range("M122").Value = -1960

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6224
$ws.Range("I122").Value = 4762
$ws.Range("J122").Value = 7686
$ws.Range("K122").Value = 14286
$ws.Range("L122").Value = 23058
$ws.Range("M122").Value = -11836
$ws.Range("N122").Value = -27958

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10710.588
$ws.Range("I7").Value = 10650.714
$ws.Range("J7").Value = 10990
$ws.Range("K7").Value = 10650.714
$ws.Range("L7").Value = 10990
$ws.Range("M7").Value = -10538.714
$ws.Range("N7").Value = -11214

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1499.75
$ws.Range("I16").Value = 1099.5
$ws.Range("K16").Value = 1099.5
$ws.Range("M16").Value = -929.5

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = ""
$ws.Range("N68").Value = ""

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = ""
$ws.Range("N71").Value = ""

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3041.5386
$ws.Range("I82").Value = 2507.1428
$ws.Range("J82").Value = 3665
$ws.Range("K82").Value = 2507.1428
$ws.Range("L82").Value = 3665
$ws.Range("M82").Value = -2146.1428
$ws.Range("N82").Value = -4387

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3041.5386
$ws.Range("I85").Value = 2507.1428
$ws.Range("J85").Value = 3665
$ws.Range("K85").Value = 2507.1428
$ws.Range("L85").Value = 3665
$ws.Range("M85").Value = -1259.1428
$ws.Range("N85").Value = -6161

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3991.818
$ws.Range("I122").Value = 3570
$ws.Range("J122").Value = 4150
$ws.Range("K122").Value = 10710
$ws.Range("L122").Value = 12450
$ws.Range("M122").Value = -8260
$ws.Range("N122").Value = -17350

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 10710.588
$ws.Range("I126").Value = 10650.714
$ws.Range("J126").Value = 10990
$ws.Range("K126").Value = 31952.142
$ws.Range("L126").Value = 32970
$ws.Range("M126").Value = -29482.142
$ws.Range("N126").Value = -37910

# WVR row 131
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970
